$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Collapse the three detailed "CORE COMPETENCIES" bullet paragraphs
#    into a single summary line.
# ---------------------------------------------------------------------
$bullet = [char]0x2022

$coreHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "CORE COMPETENCIES") {
        $coreHeading = $p
        break
    }
}

$surveyPara = $coreHeading.Next()
$redistrictingPara = $surveyPara.Next()
$dataPara = $redistrictingPara.Next()

# Remove the 2nd and 3rd detail paragraphs entirely (their ranges include
# the trailing paragraph mark, so deleting them merges the text back into
# the remaining paragraph cleanly).
$dataPara.Range.Delete()
$redistrictingPara.Range.Delete()

# Replace the remaining (first) paragraph's text with the condensed line.
$surveyPara.Range.Text = "Survey Methodology & Research Design " + $bullet + " Redistricting & Geospatial Analysis " + $bullet + " Data Analysis & Visualization"

# ---------------------------------------------------------------------
# 2. Append a new "TECHNICAL SKILLS" section at the end of the document.
#    Build all four new paragraphs in one InsertAfter call (joined by
#    carriage returns) so no stray paragraph-insertion bookkeeping is
#    left behind; then promote just the heading line to Heading2.
# ---------------------------------------------------------------------
$cr = [char]13

$skillsHeading = "TECHNICAL SKILLS"
$surveyLine = "SURVEY METHODOLOGY & RESEARCH DESIGN Survey Design and Questionnaire Development for Political and Market Research; Sampling Methodology and Statistical Analysis (R, SPSS, Stata, OSCAR); Random Device Engagement (RDE), Text Message, Web Panel, and Live Telephone Calling; Focus Groups and Qualitative Research Methodologies; Meta-analytical Dataset Development for Longitudinal Analysis; Survey Instrument Standardization and Call Methods Optimization; Expert Testimony and Consultation on Research Methodology"
$redistrictingLine = "REDISTRICTING & GEOSPATIAL ANALYSIS Redistricting Software Development and Boundary Estimation Systems; Geospatial Analysis; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Spatial Clustering and Boundary Estimation without ML Requirements; Census Data Integration and Demographic Mapping; Court Case Analysis and Expert Testimony for Redistricting; Multi-tenant Data Warehouse Design for Electoral Analytics"
$dataLine = "DATA ANALYSIS & VISUALIZATION Advanced Statistical Modeling and Analysis (Regression, Clustering, Segmentation); Data Visualization; Consumer Behavior Analysis and Market Segmentation; Machine Learning and Predictive Modeling for Targeting; Big Data Analytics; Fraud Detection and Entity Resolution Systems; Multi-million Dollar Research Project Management"

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertAfter($cr + $skillsHeading + $cr + $surveyLine + $cr + $redistrictingLine + $cr + $dataLine)

$n = $d.Paragraphs.Count
$headingPara = $d.Paragraphs.Item($n - 3)
$headingPara.Style = "Heading2"

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
